# edit.ps1 — apply the "bug in course intro" commit to the presentation.
#
# Changes:
#   1. The cached text of the (auto) date placeholder field ("1/18/24")
#      is refreshed to "3/18/24" everywhere it appears: the slide master,
#      every slide layout, and the notes master.
#   2. Slide 20 ("Homework Grades") has two small wording fixes in the
#      "Deadlines" bullet list:
#        - "request a 10 day week extension" -> "request a 10 day extension"
#        - "generous three week window"      -> "generous 10 day window"

$p = $ppt.ActivePresentation

function Update-DateShape($shape) {
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "1/18/24") {
            $tr.Text = "3/18/24"
        }
    }
}

# --- 1a. Slide master ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# --- 1b. Every slide layout ---
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# --- 1c. Notes master ---
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    Update-DateShape $notesMaster.Shapes.Item($i)
}

# --- 2. Slide 20 wording fixes ---
$slide20 = $p.Slides.Item(20)
$body = $slide20.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$txt = $tr.Text
$txt = $txt.Replace("request a 10 day week extension", "request a 10 day extension")
$txt = $txt.Replace("already generous three week window", "already generous 10 day window")
$tr.Text = $txt
